# Applies the odds updates from the 2025-04-14 FlashScore refresh.
# Each block updates one match (row) on the active sheet, identified by row number;
# column letters map to the header row (Odd_H_FT, Odd_Over05_FT, Odd_CS_x-y, etc.).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bournemouth - Fulham
$ws.Cells.Item(2, 10).Value = 1.04  # J2 (Odd_Over05_FT): 1.03 -> 1.04
$ws.Cells.Item(2, 12).Value = 1.22  # L2 (Odd_Over15_FT): 1.19 -> 1.22

# Row 4: Atl. Madrid - Valladolid
$ws.Cells.Item(4, 10).Value = 1.02  # J4 (Odd_Over05_FT): 1.03 -> 1.02
$ws.Cells.Item(4, 11).Value = 19  # K4 (Odd_Under05_FT): 17 -> 19

# Row 5: Aldosivi - Racing Club
$ws.Cells.Item(5, 11).Value = 10  # K5 (Odd_Under05_FT): 9.5 -> 10
$ws.Cells.Item(5, 14).Value = 2.08  # N5 (Odd_Over25_FT): 2.05 -> 2.08
$ws.Cells.Item(5, 15).Value = 1.73  # O5 (Odd_Under25_FT): 1.75 -> 1.73

# Row 6: Central Cordoba - Huracan
$ws.Cells.Item(6, 9).Value = 3  # I6 (Odd_A_FT): 2.9 -> 3

# Row 7: Defensa y Justicia - Union de Santa Fe
$ws.Cells.Item(7, 11).Value = 8  # K7 (Odd_Under05_FT): 7.5 -> 8

# Row 8: Godoy Cruz - Lanus
$ws.Cells.Item(8, 7).Value = 3.5  # G8 (Odd_H_FT): 3.4 -> 3.5
$ws.Cells.Item(8, 9).Value = 2.2  # I8 (Odd_A_FT): 2.25 -> 2.2
$ws.Cells.Item(8, 27).Value = 6.5  # AA8 (Odd_CS_1-1): 6 -> 6.5
$ws.Cells.Item(8, 28).Value = 21  # AB8 (Odd_CS_2-2): 19 -> 21
$ws.Cells.Item(8, 31).Value = 5.5  # AE8 (Odd_CS_0-1): 6 -> 5.5
$ws.Cells.Item(8, 32).Value = 9  # AF8 (Odd_CS_0-2): 9.5 -> 9
$ws.Cells.Item(8, 35).Value = 21  # AI8 (Odd_CS_1-3): 23 -> 21

# Row 10: Arsenal Sarandi - Patronato
$ws.Cells.Item(10, 7).Value = 2.75  # G10 (Odd_H_FT): 2.7 -> 2.75
$ws.Cells.Item(10, 8).Value = 2.8  # H10 (Odd_D_FT): 2.75 -> 2.8
$ws.Cells.Item(10, 9).Value = 2.9  # I10 (Odd_A_FT): 2.85 -> 2.9

# Row 11: Colon Santa Fe - San Telmo
$ws.Cells.Item(11, 15).Value = 1.53  # O11 (Odd_Under25_FT): 1.5 -> 1.53

# Row 12: Atletico Atlanta - Almagro
$ws.Cells.Item(12, 15).Value = 1.36  # O12 (Odd_Under25_FT): 1.33 -> 1.36

# Row 13: Nacional Potosi - Guabira
$ws.Cells.Item(13, 10).Value = 1.02  # J13 (Odd_Over05_FT): 1.03 -> 1.02
$ws.Cells.Item(13, 12).Value = 1.11  # L13 (Odd_Over15_FT): 1.14 -> 1.11

# Row 14: Londrina - Ypiranga FC
$ws.Cells.Item(14, 7).Value = 2.1  # G14 (Odd_H_FT): 2 -> 2.1
$ws.Cells.Item(14, 8).Value = 2.7  # H14 (Odd_D_FT): 2.75 -> 2.7
$ws.Cells.Item(14, 9).Value = 4.25  # I14 (Odd_A_FT): 4.6 -> 4.25
$ws.Cells.Item(14, 11).Value = 5.6  # K14 (Odd_Under05_FT): 5.7 -> 5.6
$ws.Cells.Item(14, 13).Value = 2.65  # M14 (Odd_Under15_FT): 2.67 -> 2.65
$ws.Cells.Item(14, 14).Value = 2.25  # N14 (Odd_Over25_FT): 2.22 -> 2.25
$ws.Cells.Item(14, 15).Value = 1.57  # O14 (Odd_Under25_FT): 1.6 -> 1.57
$ws.Cells.Item(14, 20).Value = 6.2  # T14 (Odd_CS_1-0): 6.1 -> 6.2
$ws.Cells.Item(14, 21).Value = 9.75  # U14 (Odd_CS_2-0): 9 -> 9.75
$ws.Cells.Item(14, 23).Value = 21  # W14 (Odd_CS_3-0): 19.5 -> 21
$ws.Cells.Item(14, 24).Value = 18.5  # X14 (Odd_CS_3-1): 17.5 -> 18.5
$ws.Cells.Item(14, 26).Value = 5.6  # Z14 (Odd_CS_0-0): 5.7 -> 5.6
$ws.Cells.Item(14, 27).Value = 5.3  # AA14 (Odd_CS_1-1): 5.4 -> 5.3
$ws.Cells.Item(14, 28).Value = 13.5  # AB14 (Odd_CS_2-2): 14 -> 13.5
$ws.Cells.Item(14, 29).Value = 70  # AC14 (Odd_CS_3-3): 75 -> 70
$ws.Cells.Item(14, 31).Value = 10  # AE14 (Odd_CS_0-1): 10.75 -> 10
$ws.Cells.Item(14, 32).Value = 24  # AF14 (Odd_CS_0-2): 27 -> 24
$ws.Cells.Item(14, 33).Value = 13.5  # AG14 (Odd_CS_1-2): 14 -> 13.5
$ws.Cells.Item(14, 34).Value = 80  # AH14 (Odd_CS_0-3): 90 -> 80
$ws.Cells.Item(14, 35).Value = 45  # AI14 (Odd_CS_1-3): 50 -> 45

# Row 15: Slavia Sofia - Botev Vratsa
$ws.Cells.Item(15, 7).Value = 1.6  # G15 (Odd_H_FT): 1.62 -> 1.6
$ws.Cells.Item(15, 8).Value = 3.8  # H15 (Odd_D_FT): 3.75 -> 3.8
$ws.Cells.Item(15, 9).Value = 6  # I15 (Odd_A_FT): 5.75 -> 6
$ws.Cells.Item(15, 10).Value = 1.05  # J15 (Odd_Over05_FT): 1.06 -> 1.05
$ws.Cells.Item(15, 11).Value = 8.5  # K15 (Odd_Under05_FT): 10 -> 8.5
$ws.Cells.Item(15, 12).Value = 1.3  # L15 (Odd_Over15_FT): 1.33 -> 1.3
$ws.Cells.Item(15, 14).Value = 2.1  # N15 (Odd_Over25_FT): 2.08 -> 2.1
$ws.Cells.Item(15, 15).Value = 1.7  # O15 (Odd_Under25_FT): 1.73 -> 1.7
$ws.Cells.Item(15, 21).Value = 6.5  # U15 (Odd_CS_2-0): 7 -> 6.5
$ws.Cells.Item(15, 26).Value = 8.5  # Z15 (Odd_CS_0-0): 9 -> 8.5
$ws.Cells.Item(15, 31).Value = 13  # AE15 (Odd_CS_0-1): 12 -> 13
$ws.Cells.Item(15, 32).Value = 29  # AF15 (Odd_CS_0-2): 26 -> 29
$ws.Cells.Item(15, 33).Value = 19  # AG15 (Odd_CS_1-2): 17 -> 19
$ws.Cells.Item(15, 34).Value = 67  # AH15 (Odd_CS_0-3): 51 -> 67
$ws.Cells.Item(15, 35).Value = 51  # AI15 (Odd_CS_1-3): 41 -> 51

# Row 17: Palestino - Union La Calera
$ws.Cells.Item(17, 14).Value = 2.03  # N17 (Odd_Over25_FT): 2 -> 2.03
$ws.Cells.Item(17, 15).Value = 1.83  # O17 (Odd_Under25_FT): 1.8 -> 1.83

# Row 20: Santa Fe - Ind. Medellin
$ws.Cells.Item(20, 7).Value = 2.55  # G20 (Odd_H_FT): 2.4 -> 2.55
$ws.Cells.Item(20, 8).Value = 3.25  # H20 (Odd_D_FT): 3.2 -> 3.25
$ws.Cells.Item(20, 9).Value = 2.7  # I20 (Odd_A_FT): 2.88 -> 2.7
$ws.Cells.Item(20, 11).Value = 8  # K20 (Odd_Under05_FT): 7.5 -> 8
$ws.Cells.Item(20, 12).Value = 1.44  # L20 (Odd_Over15_FT): 1.4 -> 1.44
$ws.Cells.Item(20, 13).Value = 2.63  # M20 (Odd_Under15_FT): 2.75 -> 2.63
$ws.Cells.Item(20, 14).Value = 2.4  # N20 (Odd_Over25_FT): 2.35 -> 2.4
$ws.Cells.Item(20, 16).Value = 1.57  # P20 (Odd_Over05_HT): -> 1.57
$ws.Cells.Item(20, 17).Value = 2.25  # Q20 (Odd_Under05_HT): 2.38 -> 2.25
$ws.Cells.Item(20, 18).Value = 2.1  # R20 (Odd_BTTS_Yes): 2 -> 2.1
$ws.Cells.Item(20, 19).Value = 1.67  # S20 (Odd_BTTS_No): 1.73 -> 1.67
$ws.Cells.Item(20, 20).Value = 6.5  # T20 (Odd_CS_1-0): 7 -> 6.5
$ws.Cells.Item(20, 22).Value = 11  # V20 (Odd_CS_2-1): 10 -> 11
$ws.Cells.Item(20, 24).Value = 26  # X20 (Odd_CS_3-1): 23 -> 26
$ws.Cells.Item(20, 26).Value = 7  # Z20 (Odd_CS_0-0): 7.5 -> 7
$ws.Cells.Item(20, 27).Value = 6.5  # AA20 (Odd_CS_1-1): 6 -> 6.5
$ws.Cells.Item(20, 28).Value = 19  # AB20 (Odd_CS_2-2): 17 -> 19
$ws.Cells.Item(20, 29).Value = 81  # AC20 (Odd_CS_3-3): 67 -> 81
$ws.Cells.Item(20, 30).Value = 1000  # AD20 (Odd_CS_4-4): 1250 -> 1000
$ws.Cells.Item(20, 31).Value = 6.5  # AE20 (Odd_CS_0-1): 7.5 -> 6.5

# Row 23: Mushuc Runa - LDU Quito
$ws.Cells.Item(23, 14).Value = 1.95  # N23 (Odd_Over25_FT): 1.98 -> 1.95
$ws.Cells.Item(23, 15).Value = 1.8  # O23 (Odd_Under25_FT): 1.83 -> 1.8

# Row 27: Bastia - Laval
$ws.Cells.Item(27, 10).Value = 1.05  # J27 (Odd_Over05_FT): 1.08 -> 1.05
$ws.Cells.Item(27, 12).Value = 1.33  # L27 (Odd_Over15_FT): 1.36 -> 1.33

# Row 28: KR Reykjavik - Valur
$ws.Cells.Item(28, 8).Value = 3.95  # H28 (Odd_D_FT): 4 -> 3.95
$ws.Cells.Item(28, 12).Value = 1.09  # L28 (Odd_Over15_FT): 1.08 -> 1.09
$ws.Cells.Item(28, 13).Value = 6.2  # M28 (Odd_Under15_FT): 6.5 -> 6.2
$ws.Cells.Item(28, 14).Value = 1.29  # N28 (Odd_Over25_FT): 1.27 -> 1.29
$ws.Cells.Item(28, 15).Value = 3.25  # O28 (Odd_Under25_FT): 3.4 -> 3.25
$ws.Cells.Item(28, 16).Value = 1.19  # P28 (Odd_Over05_HT): 1.18 -> 1.19
$ws.Cells.Item(28, 17).Value = 4.15  # Q28 (Odd_Under05_HT): 4.3 -> 4.15
$ws.Cells.Item(28, 18).Value = 1.29  # R28 (Odd_BTTS_Yes): 1.28 -> 1.29
$ws.Cells.Item(28, 19).Value = 3.3  # S28 (Odd_BTTS_No): 3.35 -> 3.3
$ws.Cells.Item(28, 20).Value = 21  # T28 (Odd_CS_1-0): 25 -> 21
$ws.Cells.Item(28, 21).Value = 27  # U28 (Odd_CS_2-0): 30 -> 27
$ws.Cells.Item(28, 22).Value = 12.5  # V28 (Odd_CS_2-1): 13 -> 12.5
$ws.Cells.Item(28, 23).Value = 50  # W28 (Odd_CS_3-0): 55 -> 50
$ws.Cells.Item(28, 25).Value = 19  # Y28 (Odd_CS_3-2): 18 -> 19
$ws.Cells.Item(28, 26).Value = 32  # Z28 (Odd_CS_0-0): 35 -> 32
$ws.Cells.Item(28, 27).Value = 9.75  # AA28 (Odd_CS_1-1): 10.25 -> 9.75
$ws.Cells.Item(28, 29).Value = 23  # AC28 (Odd_CS_3-3): 22 -> 23
$ws.Cells.Item(28, 30).Value = 90  # AD28 (Odd_CS_4-4): 80 -> 90
$ws.Cells.Item(28, 32).Value = 18  # AF28 (Odd_CS_0-2): -> 18
$ws.Cells.Item(28, 34).Value = 24  # AH28 (Odd_CS_0-3): 23 -> 24
$ws.Cells.Item(28, 36).Value = 14  # AJ28 (Odd_CS_2-3): 14.5 -> 14

# Row 30: Shamrock Rovers - Cork City
$ws.Cells.Item(30, 14).Value = 1.94  # N30 (Odd_Over25_FT): 2 -> 1.94
$ws.Cells.Item(30, 15).Value = 1.79  # O30 (Odd_Under25_FT): 1.85 -> 1.79

# Row 31: Netanya - Hapoel Haifa
$ws.Cells.Item(31, 10).Value = 1.05  # J31 (Odd_Over05_FT): 1.03 -> 1.05
$ws.Cells.Item(31, 12).Value = 1.25  # L31 (Odd_Over15_FT): 1.22 -> 1.25
$ws.Cells.Item(31, 20).Value = 9  # T31 (Odd_CS_1-0): 8.5 -> 9
$ws.Cells.Item(31, 21).Value = 12  # U31 (Odd_CS_2-0): 11 -> 12
$ws.Cells.Item(31, 24).Value = 19  # X31 (Odd_CS_3-1): 17 -> 19
$ws.Cells.Item(31, 36).Value = 29  # AJ31 (Odd_CS_2-3): 34 -> 29

# Row 32: H. Beer Sheva - Beitar Jerusalem
$ws.Cells.Item(32, 10).Value = 1.04  # J32 (Odd_Over05_FT): 1.03 -> 1.04
$ws.Cells.Item(32, 12).Value = 1.2  # L32 (Odd_Over15_FT): 1.17 -> 1.2

# Row 33: Maccabi Tel Aviv - Maccabi Haifa
$ws.Cells.Item(33, 10).Value = 1.02  # J33 (Odd_Over05_FT): 1.01 -> 1.02
$ws.Cells.Item(33, 12).Value = 1.14  # L33 (Odd_Over15_FT): 1.11 -> 1.14

# Row 34: Carpi - Legnago Salus
$ws.Cells.Item(34, 20).Value = 6.1  # T34 (Odd_CS_1-0): 6.2 -> 6.1
$ws.Cells.Item(34, 21).Value = 8.25  # U34 (Odd_CS_2-0): 8.5 -> 8.25
$ws.Cells.Item(34, 22).Value = 8.75  # V34 (Odd_CS_2-1): 8.5 -> 8.75
$ws.Cells.Item(34, 24).Value = 17.5  # X34 (Odd_CS_3-1): 17 -> 17.5
$ws.Cells.Item(34, 31).Value = 9.75  # AE34 (Odd_CS_0-1): 9.5 -> 9.75
$ws.Cells.Item(34, 33).Value = 13  # AG34 (Odd_CS_1-2): 13.5 -> 13

# Row 35: Crotone - Foggia
$ws.Cells.Item(35, 7).Value = 1.38  # G35 (Odd_H_FT): 1.36 -> 1.38
$ws.Cells.Item(35, 8).Value = 4.15  # H35 (Odd_D_FT): 4.3 -> 4.15
$ws.Cells.Item(35, 9).Value = 7.5  # I35 (Odd_A_FT): 7.4 -> 7.5
$ws.Cells.Item(35, 10).Value = 1.06  # J35 (Odd_Over05_FT): 1.05 -> 1.06
$ws.Cells.Item(35, 11).Value = 7.2  # K35 (Odd_Under05_FT): 7.7 -> 7.2
$ws.Cells.Item(35, 12).Value = 1.32  # L35 (Odd_Over15_FT): 1.28 -> 1.32
$ws.Cells.Item(35, 13).Value = 3.1  # M35 (Odd_Under15_FT): 3.35 -> 3.1
$ws.Cells.Item(35, 14).Value = 1.93  # N35 (Odd_Over25_FT): 1.83 -> 1.93
$ws.Cells.Item(35, 15).Value = 1.78  # O35 (Odd_Under25_FT): 1.88 -> 1.78
$ws.Cells.Item(35, 16).Value = 1.42  # P35 (Odd_Over05_HT): 1.39 -> 1.42
$ws.Cells.Item(35, 17).Value = 2.65  # Q35 (Odd_Under05_HT): 2.77 -> 2.65
$ws.Cells.Item(35, 18).Value = 2.27  # R35 (Odd_BTTS_Yes): 2.18 -> 2.27
$ws.Cells.Item(35, 19).Value = 1.57  # S35 (Odd_BTTS_No): 1.62 -> 1.57
$ws.Cells.Item(35, 20).Value = 5.6  # T35 (Odd_CS_1-0): 6 -> 5.6
$ws.Cells.Item(35, 21).Value = 5.7  # U35 (Odd_CS_2-0): 5.8 -> 5.7
$ws.Cells.Item(35, 24).Value = 13  # X35 (Odd_CS_3-1): 12 -> 13
$ws.Cells.Item(35, 25).Value = 37  # Y35 (Odd_CS_3-2): 35 -> 37
$ws.Cells.Item(35, 26).Value = 7.2  # Z35 (Odd_CS_0-0): 7.7 -> 7.2
$ws.Cells.Item(35, 27).Value = 8.5  # AA35 (Odd_CS_1-1): 8.75 -> 8.5
$ws.Cells.Item(35, 28).Value = 26  # AB35 (Odd_CS_2-2): 24 -> 26
$ws.Cells.Item(35, 31).Value = 16  # AE35 (Odd_CS_0-1): 17.5 -> 16
$ws.Cells.Item(35, 33).Value = 26  # AG35 (Odd_CS_1-2): 25 -> 26
$ws.Cells.Item(35, 35).Value = 110  # AI35 (Odd_CS_1-3): 100 -> 110
$ws.Cells.Item(35, 36).Value = 110  # AJ35 (Odd_CS_2-3): 100 -> 110

# Row 37: Moghreb Tetouan - IR Tanger
$ws.Cells.Item(37, 7).Value = 2.88  # G37 (Odd_H_FT): 2.9 -> 2.88
$ws.Cells.Item(37, 8).Value = 3.25  # H37 (Odd_D_FT): 3.2 -> 3.25
$ws.Cells.Item(37, 14).Value = 1.98  # N37 (Odd_Over25_FT): 2.05 -> 1.98
$ws.Cells.Item(37, 15).Value = 1.83  # O37 (Odd_Under25_FT): 1.75 -> 1.83
$ws.Cells.Item(37, 20).Value = 9.5  # T37 (Odd_CS_1-0): 9 -> 9.5
$ws.Cells.Item(37, 24).Value = 23  # X37 (Odd_CS_3-1): 26 -> 23
$ws.Cells.Item(37, 26).Value = 10  # Z37 (Odd_CS_0-0): 9 -> 10
$ws.Cells.Item(37, 30).Value = 201  # AD37 (Odd_CS_4-4): 301 -> 201
$ws.Cells.Item(37, 31).Value = 8  # AE37 (Odd_CS_0-1): 7.5 -> 8

# Row 38: Recoleta - Libertad Asuncion
$ws.Cells.Item(38, 9).Value = 1.8  # I38 (Odd_A_FT): 1.75 -> 1.8
$ws.Cells.Item(38, 11).Value = 9  # K38 (Odd_Under05_FT): 8.5 -> 9
$ws.Cells.Item(38, 26).Value = 9  # Z38 (Odd_CS_0-0): 8.5 -> 9
$ws.Cells.Item(38, 28).Value = 17  # AB38 (Odd_CS_2-2): 19 -> 17
$ws.Cells.Item(38, 29).Value = 51  # AC38 (Odd_CS_3-3): 67 -> 51
$ws.Cells.Item(38, 31).Value = 6.5  # AE38 (Odd_CS_0-1): 6 -> 6.5
$ws.Cells.Item(38, 32).Value = 8  # AF38 (Odd_CS_0-2): 7.5 -> 8

# Row 39: Deportivo Garcilaso - Alianza Lima
$ws.Cells.Item(39, 7).Value = 2.4  # G39 (Odd_H_FT): 2.38 -> 2.4
$ws.Cells.Item(39, 9).Value = 3.2  # I39 (Odd_A_FT): 3.25 -> 3.2
$ws.Cells.Item(39, 21).Value = 13  # U39 (Odd_CS_2-0): 12 -> 13
$ws.Cells.Item(39, 23).Value = 23  # W39 (Odd_CS_3-0): 21 -> 23
$ws.Cells.Item(39, 24).Value = 19  # X39 (Odd_CS_3-1): 17 -> 19
$ws.Cells.Item(39, 26).Value = 11  # Z39 (Odd_CS_0-0): 10 -> 11

# Row 43: Mafra - Penafiel
$ws.Cells.Item(43, 7).Value = 2.45  # G43 (Odd_H_FT): 2.5 -> 2.45
$ws.Cells.Item(43, 21).Value = 11  # U43 (Odd_CS_2-0): 12 -> 11

# Row 44: FC Botosani - Otelul
$ws.Cells.Item(44, 10).Value = 1.1  # J44 (Odd_Over05_FT): 1.11 -> 1.1
$ws.Cells.Item(44, 11).Value = 7  # K44 (Odd_Under05_FT): 6.5 -> 7

# Row 45: CFR Cluj - FC Rapid Bucuresti
$ws.Cells.Item(45, 7).Value = 1.6  # G45 (Odd_H_FT): 1.53 -> 1.6
$ws.Cells.Item(45, 8).Value = 3.8  # H45 (Odd_D_FT): 4 -> 3.8
$ws.Cells.Item(45, 9).Value = 5.5  # I45 (Odd_A_FT): 6.25 -> 5.5
$ws.Cells.Item(45, 18).Value = 2  # R45 (Odd_BTTS_Yes): 2.1 -> 2
$ws.Cells.Item(45, 19).Value = 1.73  # S45 (Odd_BTTS_No): 1.67 -> 1.73
$ws.Cells.Item(45, 20).Value = 6.5  # T45 (Odd_CS_1-0): 6 -> 6.5
$ws.Cells.Item(45, 21).Value = 7  # U45 (Odd_CS_2-0): 6.5 -> 7
$ws.Cells.Item(45, 28).Value = 19  # AB45 (Odd_CS_2-2): 21 -> 19
$ws.Cells.Item(45, 31).Value = 13  # AE45 (Odd_CS_0-1): 15 -> 13
$ws.Cells.Item(45, 33).Value = 17  # AG45 (Odd_CS_1-2): 19 -> 17
$ws.Cells.Item(45, 34).Value = 51  # AH45 (Odd_CS_0-3): 67 -> 51
$ws.Cells.Item(45, 35).Value = 41  # AI45 (Odd_CS_1-3): 51 -> 41

# Row 46: Napredak - Tekstilac Odzaci
$ws.Cells.Item(46, 7).Value = 1.78  # G46 (Odd_H_FT): 1.7 -> 1.78
$ws.Cells.Item(46, 8).Value = 3.45  # H46 (Odd_D_FT): 3.5 -> 3.45
$ws.Cells.Item(46, 9).Value = 4.05  # I46 (Odd_A_FT): 4.45 -> 4.05
$ws.Cells.Item(46, 11).Value = 7  # K46 (Odd_Under05_FT): 7.1 -> 7
$ws.Cells.Item(46, 18).Value = 1.87  # R46 (Odd_BTTS_Yes): 1.9 -> 1.87
$ws.Cells.Item(46, 19).Value = 1.83  # S46 (Odd_BTTS_No): 1.8 -> 1.83
$ws.Cells.Item(46, 20).Value = 6.6  # T46 (Odd_CS_1-0): 6.3 -> 6.6
$ws.Cells.Item(46, 21).Value = 8  # U46 (Odd_CS_2-0): 7.6 -> 8
$ws.Cells.Item(46, 23).Value = 14.5  # W46 (Odd_CS_3-0): 13 -> 14.5
$ws.Cells.Item(46, 26).Value = 7  # Z46 (Odd_CS_0-0): 7.1 -> 7
$ws.Cells.Item(46, 27).Value = 6.8  # AA46 (Odd_CS_1-1): 6.9 -> 6.8
$ws.Cells.Item(46, 28).Value = 16  # AB46 (Odd_CS_2-2): 16.5 -> 16
$ws.Cells.Item(46, 31).Value = 11  # AE46 (Odd_CS_0-1): 11.75 -> 11
$ws.Cells.Item(46, 32).Value = 22  # AF46 (Odd_CS_0-2): 25 -> 22
$ws.Cells.Item(46, 33).Value = 14  # AG46 (Odd_CS_1-2): 15 -> 14
$ws.Cells.Item(46, 34).Value = 65  # AH46 (Odd_CS_0-3): 75 -> 65
$ws.Cells.Item(46, 35).Value = 40  # AI46 (Odd_CS_1-3): 45 -> 40

# Row 47: Zeleznicar Pancevo - Jedinstvo U.
$ws.Cells.Item(47, 7).Value = 1.5  # G47 (Odd_H_FT): 1.52 -> 1.5
$ws.Cells.Item(47, 8).Value = 4  # H47 (Odd_D_FT): 3.9 -> 4
$ws.Cells.Item(47, 9).Value = 5.6  # I47 (Odd_A_FT): 5.4 -> 5.6
$ws.Cells.Item(47, 10).Value = 1.05  # J47 (Odd_Over05_FT): 1.06 -> 1.05
$ws.Cells.Item(47, 11).Value = 7.6  # K47 (Odd_Under05_FT): 7.5 -> 7.6
$ws.Cells.Item(47, 12).Value = 1.29  # L47 (Odd_Over15_FT): 1.3 -> 1.29
$ws.Cells.Item(47, 13).Value = 3.3  # M47 (Odd_Under15_FT): 3.25 -> 3.3
$ws.Cells.Item(47, 14).Value = 1.85  # N47 (Odd_Over25_FT): 1.88 -> 1.85
$ws.Cells.Item(47, 15).Value = 1.85  # O47 (Odd_Under25_FT): 1.83 -> 1.85
$ws.Cells.Item(47, 16).Value = 1.39  # P47 (Odd_Over05_HT): 1.4 -> 1.39
$ws.Cells.Item(47, 17).Value = 2.75  # Q47 (Odd_Under05_HT): 2.72 -> 2.75
$ws.Cells.Item(47, 21).Value = 6.6  # U47 (Odd_CS_2-0): 6.7 -> 6.6
$ws.Cells.Item(47, 23).Value = 10  # W47 (Odd_CS_3-0): 10.25 -> 10
$ws.Cells.Item(47, 24).Value = 12.5  # X47 (Odd_CS_3-1): 13 -> 12.5
$ws.Cells.Item(47, 26).Value = 7.6  # Z47 (Odd_CS_0-0): 7.5 -> 7.6
$ws.Cells.Item(47, 27).Value = 8  # AA47 (Odd_CS_1-1): 7.8 -> 8
$ws.Cells.Item(47, 31).Value = 14.5  # AE47 (Odd_CS_0-1): 13.5 -> 14.5
$ws.Cells.Item(47, 32).Value = 35  # AF47 (Odd_CS_0-2): 32 -> 35
$ws.Cells.Item(47, 33).Value = 19  # AG47 (Odd_CS_1-2): 18 -> 19
$ws.Cells.Item(47, 35).Value = 65  # AI47 (Odd_CS_1-3): 60 -> 65

# Row 49: Huesca - Malaga
$ws.Cells.Item(49, 10).Value = 1.11  # J49 (Odd_Over05_FT): 1.1 -> 1.11
$ws.Cells.Item(49, 11).Value = 6.5  # K49 (Odd_Under05_FT): 7 -> 6.5
$ws.Cells.Item(49, 12).Value = 1.53  # L49 (Odd_Over15_FT): 1.5 -> 1.53
$ws.Cells.Item(49, 13).Value = 2.38  # M49 (Odd_Under15_FT): 2.5 -> 2.38
$ws.Cells.Item(49, 14).Value = 2.7  # N49 (Odd_Over25_FT): 2.6 -> 2.7

# Row 50: Halmstad - Oster
$ws.Cells.Item(50, 11).Value = 10  # K50 (Odd_Under05_FT): 9.5 -> 10
$ws.Cells.Item(50, 12).Value = 1.33  # L50 (Odd_Over15_FT): 1.3 -> 1.33
$ws.Cells.Item(50, 13).Value = 3.25  # M50 (Odd_Under15_FT): 3.4 -> 3.25
$ws.Cells.Item(50, 15).Value = 1.75  # O50 (Odd_Under25_FT): 1.72 -> 1.75

# Row 51: AIK - Malmo FF
$ws.Cells.Item(51, 15).Value = 1.65  # O51 (Odd_Under25_FT): 1.62 -> 1.65

# Row 52: Eyupspor - Adana Demirspor
$ws.Cells.Item(52, 12).Value = 1.1  # L52 (Odd_Over15_FT): 1.07 -> 1.1
$ws.Cells.Item(52, 14).Value = 1.33  # N52 (Odd_Over25_FT): 1.3 -> 1.33
$ws.Cells.Item(52, 18).Value = 1.87  # R52 (Odd_BTTS_Yes): 1.91 -> 1.87
$ws.Cells.Item(52, 19).Value = 1.77  # S52 (Odd_BTTS_No): 1.8 -> 1.77

# Row 54: Racing Montevideo - Liverpool M.
$ws.Cells.Item(54, 10).Value = 1.08  # J54 (Odd_Over05_FT): 1.05 -> 1.08
$ws.Cells.Item(54, 12).Value = 1.44  # L54 (Odd_Over15_FT): 1.41 -> 1.44
$ws.Cells.Item(54, 13).Value = 2.63  # M54 (Odd_Under15_FT): 2.62 -> 2.63
$ws.Cells.Item(54, 15).Value = 1.53  # O54 (Odd_Under25_FT): 1.5 -> 1.53

# Row 55: Defensor Sp. - CA Cerro
$ws.Cells.Item(55, 10).Value = 1.1  # J55 (Odd_Over05_FT): 1.07 -> 1.1
$ws.Cells.Item(55, 12).Value = 1.44  # L55 (Odd_Over15_FT): 1.41 -> 1.44
$ws.Cells.Item(55, 13).Value = 2.63  # M55 (Odd_Under15_FT): 2.62 -> 2.63
$ws.Cells.Item(55, 15).Value = 1.57  # O55 (Odd_Under25_FT): 1.54 -> 1.57
